# Apply edits to the "Score" worksheet of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")
$ws.Activate()

# New values for column O (quiz/assignment score that drives P=O/2, Q=SUM(F:N,P), R=VLOOKUP(Q,...))
$oValues = @{
    11 = 28
    12 = 30
    15 = 34
    16 = 30
    17 = 26
    18 = 39
    20 = 30
    21 = 31
    22 = 39
    23 = 33
    26 = 19
    27 = 30
    29 = 41
}

foreach ($row in $oValues.Keys) {
    $ws.Range("O$row").Value = $oValues[$row]
}

# Update the view: top-left visible cell (A7) and active selection (O27).
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("O27").Select()

# Update print scale for the sheet's page setup (manual scale, keeps "fit to
# height" unlimited as in the source file).
$ws.PageSetup.Zoom = 80
$ws.PageSetup.FitToPagesTall = $false
